$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D/E cells in this sheet are stored as text (e.g. "67.558.71", "0.998",
# "  +0.96%  "). Force text format before assigning so Excel does not
# auto-convert plain-decimal-looking values (e.g. "0.999") into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.665.17"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.873.25"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "461.98"
$ws.Range("E5").Value = "  +9.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.13"
$ws.Range("E6").Value = "  +14.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.751"
$ws.Range("E9").Value = "  +4.41%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("E11").Value = "  -4.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.10"
$ws.Range("E12").Value = "  +8.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.46"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.479.85"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.85"
$ws.Range("E15").Value = "  -4.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.860.08"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.15"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +7.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.751.21"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.24"
$ws.Range("E21").Value = "  +4.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.92"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.29"
$ws.Range("E23").Value = "  +8.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.56"
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.56"
$ws.Range("E25").Value = "  +9.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.46"
$ws.Range("E26").Value = "  +11.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.78"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "754.40"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  +11.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.78"
$ws.Range("E32").Value = "  +6.00%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.63"
$ws.Range("E34").Value = "  +13.38%  "
$ws.Range("E35").Value = "  +7.67%  "
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  +5.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.355"
$ws.Range("E40").Value = "  +12.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +16.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0682"
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("E44").Value = "  +5.41%  "
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.46"
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.27"
$ws.Range("E47").Value = "  +5.20%  "
$ws.Range("E48").Value = "  +7.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.14"
$ws.Range("E49").Value = "  +5.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.09"
$ws.Range("E50").Value = "  +3.10%  "
$ws.Range("E51").Value = "  +3.15%  "
